$d = $word.ActiveDocument

# --- Create the "Normal (Web)" style (mirrors what Word mints the first time the style is used) ---
$normalWeb = $d.Styles.Add("Normal (Web)", 1)
$normalWeb.BaseStyle = "Normal"
$normalWeb.Priority = 99
$normalWeb.UnhideWhenUsed = $true
$normalWeb.Font.NameAscii = "Times New Roman"
$normalWeb.Font.NameFarEast = "Times New Roman"
$normalWeb.Font.NameOther = "Times New Roman"
$normalWeb.Font.NameBi = "Times New Roman"
$normalWeb.Font.Name = "Times New Roman"
$normalWeb.Font.Size = 12
$normalWeb.Font.SizeBi = 12
$normalWeb.ParagraphFormat.SpaceBefore = 5
$normalWeb.ParagraphFormat.SpaceBeforeAuto = $true
$normalWeb.ParagraphFormat.SpaceAfter = 5
$normalWeb.ParagraphFormat.SpaceAfterAuto = $true
$normalWeb.ParagraphFormat.LineSpacingRule = 0

# --- Locate the "6.1 Fedora vs Ubuntu" heading paragraph ---
$rng = $d.Content
$null = $rng.Find.Execute("Fedora vs Ubuntu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# --- Insert the six new body paragraphs after the heading, each styled Normal (Web) ---
$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Fedora Linux is a Linux distribution which has been developed by the Fedora Project []. There are several well know Linux distributions which are broadly used throughout industry and academia, each of which poses their own unique characteristics and user experience. Fedora is one of the most popular Linux based distributions. It contains software distributed under various sources, including free and open-sources, and is currently identified as a distribution which is a leader in open-source technologies [2]. Fedora is also a fundamental source for Red Hat Enterprise Linux and every six months a new version of Fedora Linux is released [3]. Additionally, Fedora can be installed and run on a desktop or laptop, but also can be installed and run on a virtual machine which ensures that the Linux distribution can be utilized by users across several types of operating systems. ")
$p0 = $rng.Paragraphs(1)
$p0.Style = "Normal (Web)"
$p0.SpaceBefore = 0
$p0.SpaceBeforeAuto = $false
$p0.SpaceAfter = 8
$p0.SpaceAfterAuto = $false

$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("The Ubuntu Linux distribution is composed of free and open-source software. It is based on the Debian Linux based distribution and as a result some of its characteristics resemble Debian [4]. Specifically, the Ubuntu architecture and infrastructure has primarily been built off Debian architecture and infrastructure. Ubuntu packages are also based on the packages of the Debian Linux distribution with both distributions utilizing deb package format as well as package management tools [4]. However, this does not mean that the packages of these two Linux distributions are compatible. Ubuntu is known for having three active editions, namely the desktop, server, and core editions. The core edition is a popular Linux distribution in the robotics community and is often used for Internet of Things (IoT) devices [4]. Other popular uses of Ubuntu include cloud computing. Furthermore, Ubuntu can be installed and run on a desktop or laptop, but also can be installed and run on a virtual machine which ensures that the Linux distribution can be utilized by users across several types of operating systems.")
$p1 = $rng.Paragraphs(1)
$p1.Style = "Normal (Web)"
$p1.SpaceBefore = 0
$p1.SpaceBeforeAuto = $false
$p1.SpaceAfter = 8
$p1.SpaceAfterAuto = $false

$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Both Fedora and Ubuntu distributions possess positive and negative characteristics that were observed while during user interaction. These differences were primarily observed at the Graphical User Interface (GUI) level of the distributions. However, while installing and setting up the distributions on a virtual machine, several characteristics were observed that are important to mention in this comparison since the installation process is the first point of contact between the user and the distribution. When installing Fedora and Ubuntu on a desktop or in a virtual machine, the entire process takes between 5 – 10 minutes [5]. This is a positive characteristic as it ensures that the wait time for the user to begin engaging with the distribution is short. Additionally, the installation process also ensures that all updates are included in the installation process to ensure that users begin with the most up to date version of the distribution [5]. This is a very positive quality that both distributions possess and is a great starting point for user interaction with these distributions as it ensures maximum chance of a positive user experience.")
$p2 = $rng.Paragraphs(1)
$p2.Style = "Normal (Web)"
$p2.SpaceBefore = 0
$p2.SpaceBeforeAuto = $false
$p2.SpaceAfter = 8
$p2.SpaceAfterAuto = $false

$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("The Graphical User Interface (GUI) of Ubuntu and Fedora both use the GNOME desktop environment by default [5]. However, Fedora uses the standard GNOME GUI whereas Ubuntu has customized the GNOME GUI so that it resembles a Unity desktop.  Upon initializing the Ubuntu and Fedora Linux distributions, several differences on the desktop are observed. Namely, the activities bar on the left-hand side of the desktop is much larger for Ubuntu than it is for Fedora. Using a larger activities bar may be desirable for users because it helps to clearly track activities across several applications at once. However, Fedora has a search bar present at the top of the desktop and Ubuntu does not have a visible search bar. This search bar is a very positive feature since it allows the user to navigate the system a lot more freely. Additionally, the fact that the search bar is almost in alignment with the users' eyes is great because it directs the user straight to the search bar to find anything they need. Another difference that was observed was that Ubuntu had the ability to minimize applications once they were opened. The method was similar to the method used on a windows operating system. However, Fedora did not have an intuitive method of minimizing applications once they were opened. Thus, making it difficult to navigate between applications while using the system.")
$p3 = $rng.Paragraphs(1)
$p3.Style = "Normal (Web)"
$p3.SpaceBefore = 0
$p3.SpaceBeforeAuto = $false
$p3.SpaceAfter = 8
$p3.SpaceAfterAuto = $false

$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("The Ubuntu distribution supports a variety of applications which helps ensure that the user is comfortable with the accessibility standpoint of the Ubuntu distribution [15]. For example, users have access to commonly used integrated development environments, communication applications for instant messaging, internet access through a variety of applications such as Firefox and Chromium, as well as music streaming applications such as Spotify. Similarly, Fedora possesses the same accessibility characteristics as Ubuntu. Both GUI’s support the ability to customize the desktop experience which helps ensure that the user experience can be tailored to the user's needs.  The terminal is easily accessible while working on Ubuntu and Fedora. This is a fundamental characteristic as the entire system can be accessed through the terminal and most Linux users require terminal access. ")
$p4 = $rng.Paragraphs(1)
$p4.Style = "Normal (Web)"
$p4.SpaceBefore = 0
$p4.SpaceBeforeAuto = $false
$p4.SpaceAfter = 8
$p4.SpaceAfterAuto = $false

$rng.InsertParagraphAfter()
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter("Other superficial differences between the Fedora and Ubuntu distributions were observed when becoming familiar with the two distributions. Some of these differences included the theme color of the desktop and the slightly different locations of the applications. While these differences do not yield a specific positive or negative trait, it is worth mentioning due to the fact that it may become more relevant at a later stage in the project.  Overall, it would appear that the Ubuntu and Fedora Linux distributions have many similarities in terms of graphical user interface and user experience. This aligns well with the fact that these two Linux distributions are among the most popular Linux distributions.")
$p5 = $rng.Paragraphs(1)
$p5.Style = "Normal (Web)"
$p5.SpaceBefore = 0
$p5.SpaceBeforeAuto = $false
$p5.SpaceAfter = 8
$p5.SpaceAfterAuto = $false

Write-Output "done"